$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $ref, $val) {
    $sheet.Range($ref).Value = "'" + $val
    $sheet.Range($ref).Style = "Normal"
}

# Row 2
Set-TextCell $ws 'D2' '71.106.97'
$ws.Range('E2').Value = '  +2.69%  '

# Row 3
Set-TextCell $ws 'D3' '3.575.80'
$ws.Range('E3').Value = '  +1.27%  '

# Row 4
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
Set-TextCell $ws 'D5' '583.32'
$ws.Range('E5').Value = '  +2.11%  '

# Row 6
Set-TextCell $ws 'D6' '186.09'
$ws.Range('E6').Value = '  +2.10%  '

# Row 7
Set-TextCell $ws 'D7' '3.566.40'
$ws.Range('E7').Value = '  +1.17%  '

# Row 8
Set-TextCell $ws 'D8' '0.622'
$ws.Range('E8').Value = '  +1.14%  '

# Row 9
$ws.Range('E9').Value = '  -0.04%  '

# Row 10
Set-TextCell $ws 'D10' '0.215'
$ws.Range('E10').Value = '  +14.78%  '

# Row 11
Set-TextCell $ws 'D11' '0.653'
$ws.Range('E11').Value = '  +2.12%  '

# Row 12
Set-TextCell $ws 'D12' '54.38'
$ws.Range('E12').Value = '  +1.43%  '

# Row 13
Set-TextCell $ws 'D13' '0.0000318'
$ws.Range('E13').Value = '  +5.70%  '

# Row 14
Set-TextCell $ws 'D14' '9.53'
$ws.Range('E14').Value = '  +0.67%  '

# Row 15
Set-TextCell $ws 'D15' '4.134.49'
$ws.Range('E15').Value = '  +0.88%  '

# Row 16
Set-TextCell $ws 'D16' '71.130.33'
$ws.Range('E16').Value = '  +2.75%  '

# Row 17
Set-TextCell $ws 'D17' '19.23'
$ws.Range('E17').Value = '  -0.28%  '

# Row 18
Set-TextCell $ws 'D18' '3.534.48'
$ws.Range('E18').Value = '  +0.74%  '

# Row 19
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws 'D19' '568.41'
$ws.Range('E19').Value = '  +5.35%  '

# Row 20
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws 'D20' '12.36'
$ws.Range('E20').Value = '  -0.04%  '

# Row 21
Set-TextCell $ws 'D21' '0.120'
$ws.Range('E21').Value = '  +0.63%  '

# Row 22
$ws.Range('E22').Value = '  -1.57%  '

# Row 23
Set-TextCell $ws 'D23' '17.64'
$ws.Range('E23').Value = '  -9.74%  '

# Row 24
$ws.Range('E24').Value = '  +3.14%  '

# Row 25
Set-TextCell $ws 'D25' '4.58'
$ws.Range('E25').Value = '  +5.36%  '

# Row 26
Set-TextCell $ws 'D26' '94.76'
$ws.Range('E26').Value = '  +0.66%  '

# Row 27
Set-TextCell $ws 'D27' '11.26'
$ws.Range('E27').Value = '  +1.52%  '

# Row 28
Set-TextCell $ws 'D28' '2.94'
$ws.Range('E28').Value = '  +1.12%  '

# Row 29
Set-TextCell $ws 'D29' '9.15'
$ws.Range('E29').Value = '  +0.76%  '

# Row 30
Set-TextCell $ws 'D30' '32.57'
$ws.Range('E30').Value = '  +2.51%  '

# Row 31
Set-TextCell $ws 'D31' '7.26'
$ws.Range('E31').Value = '  -1.21%  '

# Row 32
Set-TextCell $ws 'D32' '12.30'
$ws.Range('E32').Value = '  -1.87%  '

# Row 33
Set-TextCell $ws 'D33' '0.115'
$ws.Range('E33').Value = '  +0.93%  '

# Row 34
Set-TextCell $ws 'D34' '64.09'
$ws.Range('E34').Value = '  -1.14%  '

# Row 35
Set-TextCell $ws 'D35' '3.35'
$ws.Range('E35').Value = '  +6.97%  '

# Row 36
Set-TextCell $ws 'D36' '549.13'
$ws.Range('E36').Value = '  -3.81%  '

# Row 37
$ws.Range('E37').Value = '  +5.15%  '

# Row 38
Set-TextCell $ws 'D38' "0.0$([char]0x2083)0803"
$ws.Range('E38').Value = '  +5.27%  '

# Row 39
Set-TextCell $ws 'D39' '37.63'
$ws.Range('E39').Value = '  -1.25%  '

# Row 40
$ws.Range('E40').Value = '  -0.02%  '

# Row 41
Set-TextCell $ws 'D41' '3.28'
$ws.Range('E41').Value = '  +6.07%  '

# Row 42
Set-TextCell $ws 'D42' '3.501.95'
$ws.Range('E42').Value = '  +10.98%  '

# Row 43
Set-TextCell $ws 'D43' '3.47'
$ws.Range('E43').Value = '  +3.25%  '

# Row 44
$ws.Range('E44').Value = '  +1.99%  '

# Row 45
Set-TextCell $ws 'D45' '0.0446'
$ws.Range('E45').Value = '  +0.93%  '

# Row 46
$ws.Range('E46').Value = '  +0.58%  '

# Row 47
$ws.Range('E47').Value = '  -0.41%  '

# Row 48
Set-TextCell $ws 'D48' '9.26'
$ws.Range('E48').Value = '  +0.61%  '

# Row 49
$ws.Range('E49').Value = '  +2.34%  '

# Row 50
$ws.Range('E50').Value = '  +0.36%  '

# Row 51
Set-TextCell $ws 'D51' '1.45'
$ws.Range('E51').Value = '  +4.18%  '

